$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing columns B..F to C..G
$ws.Columns("B:B").Insert()

# The inserted column inherits formatting from column A for the data rows;
# the original layout had no style on those cells, so clear it.
$ws.Range("B2:B20").ClearFormats()

# Set header for new column B ("segments") and give it the same look as the
# other header cells (bold, centered, thin border) like the rest of row 1,
# by copying the formatting from the neighboring header cell.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the segment names from column A into the new column B, and put a
# numeric index (0-based) into column A instead.
$names = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 1).Value = $i
}
